# Feat: Atualização da Classificação de Defeitos
#
# Renames the catalog header cells (A1/B1) to reflect the new
# "Classificação de Defeitos" wording, widens column A to fit the new
# failure-code header, normalises the header rows back to the sheet's
# default row height, and moves the active selection to D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the header text (CÓDIGO -> CÓDIGO DA FALHA, DESCRIÇÃO DO
#    MATERIAL -> DESCRIÇÃO DA FALHA).
$ws.Range("A1").Value = "CÓDIGO DA FALHA"
$ws.Range("B1").Value = "DESCRIÇÃO DA FALHA"

# 2) Give column A an explicit width so the new header fits.
$ws.Columns.Item(1).ColumnWidth = 16.6

# 3) Rows 1 and 2 previously carried an explicit (smaller) row height;
#    auto-fit them back so they follow the sheet's default height again.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()

# 4) Move the active selection to D10.
$ws.Range("D10").Select() | Out-Null
